# Restored from revision of admin on 04/30/2021 12:06:25 PM.TEST Author: admin. Type: SAVE.
# The only functional change: cell C10 ("Integer min" for rule R20) changes
# from 18 to 1. Cell style/formatting (s="20") is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
